# Update PFT (column B) values for specific species rows
# per commit: "add functions for pft from family and woodiness"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    13  = "not found"
    48  = "not found"
    55  = "not found"
    56  = "not assigned"
    72  = "not found"
    76  = "not found"
    78  = "not found"
    101 = "not found"
    103 = "not found"
    125 = "not found"
    128 = "not found"
    171 = "not found"
    174 = "not found"
    207 = "not found"
    221 = "not found"
    229 = "not found"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
